$d = $word.ActiveDocument

# Target paragraph is the 2nd paragraph in the document body (the contract intro
# paragraph containing the "…, с одной стороны, и  FORMTEXT  именуемое в дальнейшем
# «Заказчик» ..." text together with the legacy form-field / bookmark markup).
$p2 = $d.Paragraphs(2)
$full = $p2.Range
# Exclude the trailing paragraph mark so only the paragraph's run content is
# replaced; the paragraph's own mark-run properties (pPr/rPr) are left as-is
# (and we restate them below to be explicit/safe).
$inner = $d.Range($full.Start, $full.End - 1)

$xmlFrag = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/dummy.xml" pkg:contentType="application/xml">
<pkg:xmlData>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00E61492" w:rsidRPr="00F12767" w:rsidRDefault="00DB2A42" w:rsidP="006B242C"><w:pPr><w:rPr><w:b /></w:rPr></w:pPr><w:proofErr w:type="gramStart" /><w:r><w:rPr><w:b /></w:rPr><w:t>…</w:t></w:r><w:r w:rsidR="00A47FB6" w:rsidRPr="00EC6BBD"><w:t xml:space="preserve">, с одной стороны, и </w:t></w:r><w:r w:rsidR="00B70A06" w:rsidRPr="00B70A06"><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="009B3357"><w:rPr><w:rFonts w:eastAsia="Times-Roman" /><w:color w:val="000000" w:themeColor="text1" /><w:lang w:val="en-US" /></w:rPr><w:fldChar w:fldCharType="begin"><w:ffData><w:name w:val="t1" /><w:enabled w:val="0" /><w:calcOnExit w:val="0" /><w:textInput><w:default w:val="Text" /></w:textInput></w:ffData></w:fldChar></w:r><w:bookmarkStart w:id="0" w:name="t1" /><w:r><w:rPr><w:rFonts w:eastAsia="Times-Roman" /><w:color w:val="000000" w:themeColor="text1" /></w:rPr><w:instrText xml:space="preserve"> </w:instrText></w:r><w:r w:rsidR="009B3357"><w:rPr><w:rFonts w:eastAsia="Times-Roman" /><w:color w:val="000000" w:themeColor="text1" /><w:lang w:val="en-US" /></w:rPr><w:instrText>FORMTEXT</w:instrText></w:r><w:r><w:rPr><w:rFonts w:eastAsia="Times-Roman" /><w:color w:val="000000" w:themeColor="text1" /></w:rPr><w:instrText xml:space="preserve"> </w:instrText></w:r><w:r w:rsidR="009B3357"><w:rPr><w:rFonts w:eastAsia="Times-Roman" /><w:color w:val="000000" w:themeColor="text1" /><w:lang w:val="en-US" /></w:rPr></w:r><w:r w:rsidR="009B3357"><w:rPr><w:rFonts w:eastAsia="Times-Roman" /><w:color w:val="000000" w:themeColor="text1" /><w:lang w:val="en-US" /></w:rPr><w:fldChar w:fldCharType="separate" /></w:r><w:r w:rsidR="009B3357"><w:rPr><w:rFonts w:eastAsia="Times-Roman" /><w:noProof /><w:color w:val="000000" w:themeColor="text1" /><w:lang w:val="en-US" /></w:rPr><w:t>Text</w:t></w:r><w:r w:rsidR="009B3357"><w:rPr><w:rFonts w:eastAsia="Times-Roman" /><w:color w:val="000000" w:themeColor="text1" /><w:lang w:val="en-US" /></w:rPr><w:fldChar w:fldCharType="end" /></w:r><w:bookmarkEnd w:id="0" /><w:r w:rsidR="00B56DC0" w:rsidRPr="00BA6C2C"><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="00A47FB6" w:rsidRPr="00DE147A"><w:rPr><w:rFonts w:eastAsia="Times-Roman" /><w:color w:val="000000" w:themeColor="text1" /></w:rPr><w:t xml:space="preserve">именуемое </w:t></w:r><w:r w:rsidR="003A1F8F"><w:rPr><w:rFonts w:eastAsia="Times-Roman" /><w:color w:val="000000" w:themeColor="text1" /></w:rPr><w:t>в</w:t></w:r><w:r w:rsidR="00206CB3"><w:rPr><w:rFonts w:eastAsia="Times-Roman" /><w:color w:val="000000" w:themeColor="text1" /></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsia="Times-Roman" /><w:color w:val="000000" w:themeColor="text1" /><w:lang w:val="en-US" /></w:rPr><w:fldChar w:fldCharType="begin"><w:ffData><w:name w:val="t2" /><w:enabled /><w:calcOnExit w:val="0" /><w:textInput><w:default w:val="Text" /></w:textInput></w:ffData></w:fldChar></w:r><w:bookmarkStart w:id="1" w:name="t2" /><w:r><w:rPr><w:rFonts w:eastAsia="Times-Roman" /><w:color w:val="000000" w:themeColor="text1" /></w:rPr><w:instrText xml:space="preserve"> </w:instrText></w:r><w:r><w:rPr><w:rFonts w:eastAsia="Times-Roman" /><w:color w:val="000000" w:themeColor="text1" /><w:lang w:val="en-US" /></w:rPr><w:instrText>FORMTEXT</w:instrText></w:r><w:r><w:rPr><w:rFonts w:eastAsia="Times-Roman" /><w:color w:val="000000" w:themeColor="text1" /></w:rPr><w:instrText xml:space="preserve"> </w:instrText></w:r><w:r><w:rPr><w:rFonts w:eastAsia="Times-Roman" /><w:color w:val="000000" w:themeColor="text1" /><w:lang w:val="en-US" /></w:rPr></w:r><w:r><w:rPr><w:rFonts w:eastAsia="Times-Roman" /><w:color w:val="000000" w:themeColor="text1" /><w:lang w:val="en-US" /></w:rPr><w:fldChar w:fldCharType="separate" /></w:r><w:r><w:rPr><w:rFonts w:eastAsia="Times-Roman" /><w:noProof /><w:color w:val="000000" w:themeColor="text1" /><w:lang w:val="en-US" /></w:rPr><w:t>Text</w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsia="Times-Roman" /><w:color w:val="000000" w:themeColor="text1" /><w:lang w:val="en-US" /></w:rPr><w:fldChar w:fldCharType="end" /></w:r><w:bookmarkEnd w:id="1" /><w:r><w:rPr><w:rFonts w:eastAsia="Times-Roman" /><w:color w:val="000000" w:themeColor="text1" /></w:rPr><w:t xml:space="preserve"> дальнейшем «Заказчик», в лице</w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsia="Times-Roman" /><w:color w:val="000000" w:themeColor="text1" /></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>${</w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsia="Times-Roman" /><w:color w:val="000000" w:themeColor="text1" /></w:rPr><w:t>SOME_TEXT</w:t></w:r><w:r><w:t>}</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsia="Times-Roman" /><w:color w:val="000000" w:themeColor="text1" /></w:rPr><w:t xml:space="preserve">действующего на основании </w:t></w:r><w:proofErr w:type="gramEnd" /><w:r><w:fldChar w:fldCharType="begin"><w:ffData><w:name w:val="t3" /><w:enabled /><w:calcOnExit w:val="0" /><w:textInput><w:default w:val="Text" /></w:textInput></w:ffData></w:fldChar></w:r><w:bookmarkStart w:id="2" w:name="t3" /><w:r><w:instrText xml:space="preserve"> FORMTEXT </w:instrText></w:r><w:r><w:fldChar w:fldCharType="separate" /></w:r><w:r><w:rPr><w:noProof /></w:rPr><w:t>Text</w:t></w:r><w:r><w:fldChar w:fldCharType="end" /></w:r><w:bookmarkEnd w:id="2" /><w:proofErr w:type="gramStart" /><w:r><w:rPr><w:rFonts w:eastAsia="Times-Roman" /><w:color w:val="000000" w:themeColor="text1" /></w:rPr><w:t>с другой стороны, заключили настоящий Договор о нижеследующем:</w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsia="Times-Roman" /><w:color w:val="000000" w:themeColor="text1" /></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:b /></w:rPr><w:t>…</w:t></w:r><w:proofErr w:type="gramEnd" /></w:p>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$inner.InsertXML($xmlFrag)
